# Scheduled-runner update: refresh computed pricing/profit figures across the
# per-job "Sheets" (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) in the workbook.
# Each block below updates the recalculated currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) for the specific leve rows that changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 365.83334
$ws.Range("I38").Value = 365.83334
$ws.Range("K38").Value = 1097.50002
$ws.Range("M38").Value = -725.5000199999999

$ws.Range("H69").Value = 25711.428
$ws.Range("I69").Value = 12998
$ws.Range("K69").Value = 38994
$ws.Range("M69").Value = -38120

$ws.Range("H72").Value = 25711.428
$ws.Range("I72").Value = 12998
$ws.Range("K72").Value = 116982
$ws.Range("M72").Value = -112614

$ws.Range("H76").Value = 4210
$ws.Range("I76").Value = 3900
$ws.Range("J76").Value = 5450
$ws.Range("K76").Value = 3900
$ws.Range("L76").Value = 5450
$ws.Range("M76").Value = -3585
$ws.Range("N76").Value = -6080

$ws.Range("H79").Value = 4210
$ws.Range("I79").Value = 3900
$ws.Range("J79").Value = 5450
$ws.Range("K79").Value = 3900
$ws.Range("L79").Value = 5450
$ws.Range("M79").Value = -2808
$ws.Range("N79").Value = -7634

$ws.Range("H98").Value = 47643570
$ws.Range("I98").Value = 52658580
$ws.Range("J98").Value = 944
$ws.Range("K98").Value = 52658580
$ws.Range("L98").Value = 944
$ws.Range("M98").Value = -52657082
$ws.Range("N98").Value = -3940

$ws.Range("H111").Value = 4836.4287
$ws.Range("I111").Value = 4976
$ws.Range("J111").Value = 3999
$ws.Range("K111").Value = 14928
$ws.Range("L111").Value = 11997
$ws.Range("M111").Value = -11861
$ws.Range("N111").Value = -18131

$ws.Range("H113").Value = 45456696
$ws.Range("I113").Value = 12502346
$ws.Range("K113").Value = 12502346
$ws.Range("M113").Value = -12499092

$ws.Range("H122").Value = 47643570
$ws.Range("I122").Value = 52658580
$ws.Range("J122").Value = 944
$ws.Range("K122").Value = 157975740
$ws.Range("L122").Value = 2832
$ws.Range("M122").Value = -157973290
$ws.Range("N122").Value = -7732

$ws.Range("H127").Value = 1636.7
$ws.Range("J127").Value = 1449.5
$ws.Range("L127").Value = 4348.5
$ws.Range("N127").Value = -14268.5

$ws.Range("H135").Value = 1522.3077
$ws.Range("I135").Value = 1522.3077
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13700.7693
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11165.7693
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1353.7693
$ws.Range("J2").Value = 1084.3334
$ws.Range("L2").Value = 1084.3334
$ws.Range("N2").Value = -1310.3334

$ws.Range("H32").Value = 10003050
$ws.Range("I32").Value = 10206378
$ws.Range("K32").Value = 10206378
$ws.Range("M32").Value = -10206091

$ws.Range("H61").Value = 6961871
$ws.Range("J61").Value = 106953.5
$ws.Range("L61").Value = 106953.5
$ws.Range("N61").Value = -107377.5

$ws.Range("H74").Value = 14715698
$ws.Range("I74").Value = 25001666
$ws.Range("K74").Value = 25001666
$ws.Range("M74").Value = -25000792

$ws.Range("H77").Value = 14715698
$ws.Range("I77").Value = 25001666
$ws.Range("K77").Value = 125008330
$ws.Range("M77").Value = -125003962

$ws.Range("H111").Value = 93250
$ws.Range("J111").Value = 93250
$ws.Range("L111").Value = 93250
$ws.Range("N111").Value = -101430

$ws.Range("H116").Value = 1353.7693
$ws.Range("J116").Value = 1084.3334
$ws.Range("L116").Value = 1084.3334
$ws.Range("N116").Value = -5672.3334

$ws.Range("H136").Value = 6961871
$ws.Range("J136").Value = 106953.5
$ws.Range("L136").Value = 320860.5
$ws.Range("N136").Value = -325960.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1353.7693
$ws.Range("J3").Value = 1084.3334
$ws.Range("L3").Value = 1084.3334
$ws.Range("N3").Value = -1312.3334

$ws.Range("H20").Value = 4498.4375
$ws.Range("I20").Value = 4154.7856
$ws.Range("J20").Value = 6904
$ws.Range("K20").Value = 4154.7856
$ws.Range("L20").Value = 6904
$ws.Range("M20").Value = -3907.7856
$ws.Range("N20").Value = -7398

$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -52122

$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -160608

$ws.Range("H99").Value = 2510.0833
$ws.Range("I99").Value = 2035.6666
$ws.Range("K99").Value = 2035.6666
$ws.Range("M99").Value = -537.6666

$ws.Range("H107").Value = 5431.5713
$ws.Range("I107").Value = 4250
$ws.Range("K107").Value = 4250
$ws.Range("M107").Value = -2330

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H134").Value = 40884
$ws.Range("I134").Value = 771.36365
$ws.Range("J134").Value = 261503.5
$ws.Range("K134").Value = 2314.09095
$ws.Range("L134").Value = 784510.5
$ws.Range("M134").Value = 220.9090500000002
$ws.Range("N134").Value = -789580.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1195.3
$ws.Range("I132").Value = 1222.6111
$ws.Range("K132").Value = 3667.8333
$ws.Range("M132").Value = -1137.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3953834.8
$ws.Range("I4").Value = 10088066
$ws.Range("J4").Value = 79583.78999999999
$ws.Range("K4").Value = 30264198
$ws.Range("L4").Value = 238751.37
$ws.Range("M4").Value = -30264086
$ws.Range("N4").Value = -238975.37

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 10000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H107").Value = 1672.5834
$ws.Range("I107").Value = 1104.375
$ws.Range("K107").Value = 1104.375
$ws.Range("M107").Value = 815.625

$ws.Range("H132").Value = 66669110
$ws.Range("I132").Value = 71431100
$ws.Range("J132").Value = 1194
$ws.Range("K132").Value = 214293300
$ws.Range("L132").Value = 3582
$ws.Range("M132").Value = -214290770
$ws.Range("N132").Value = -8642

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 18498.334
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 18498.334
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 18498.334
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -19820.334

$ws.Range("H132").Value = 345454.84
$ws.Range("I132").Value = 346605.2
$ws.Range("K132").Value = 1039815.6
$ws.Range("M132").Value = -1037285.6

$ws.Range("H136").Value = 41217.758
$ws.Range("I136").Value = 6821.9443
$ws.Range("K136").Value = 20465.8329
$ws.Range("M136").Value = -17915.8329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 169833
$ws.Range("J5").Value = 4749.5
$ws.Range("L5").Value = 4749.5
$ws.Range("N5").Value = -4973.5
